{"js": "// Remove the task paragraph \"\u041e\u0431\u043d\u043e\u0432\u0438\u0442\u044c \u0431\u0434 \u2013 \u0434\u043e\u0431\u0430\u0432\u0438\u0442\u044c \u0432\u043e\u0437\u0440\u0430\u0441\u0442\" entirely.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(p => p.text.indexOf(\"\u041e\u0431\u043d\u043e\u0432\u0438\u0442\u044c\") !== -1 && p.text.indexOf(\"\u0434\u043e\u0431\u0430\u0432\u0438\u0442\u044c \u0432\u043e\u0437\u0440\u0430\u0441\u0442\") !== -1);\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*\u041e\u0431\u043d\u043e\u0432\u0438\u0442\u044c*\" -and $t -like \"*\u0434\u043e\u0431\u0430\u0432\u0438\u0442\u044c \u0432\u043e\u0437\u0440\u0430\u0441\u0442*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
